$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D8").Value = 2417100
$ws.Range("E8").Value = 2499600
$ws.Range("F8").Value = 2192400
$ws.Range("G8").Value = 2530000
$ws.Range("H8").Value = 2359800
$ws.Range("I8").Value = 2704700
$ws.Range("J8").Value = 2303500

$ws.Range("D9").Value = 1293300
$ws.Range("E9").Value = 1266500
$ws.Range("F9").Value = 1070300
$ws.Range("G9").Value = 1178200
$ws.Range("H9").Value = 1082100
$ws.Range("I9").Value = 1393400
$ws.Range("J9").Value = 1022600

$ws.Range("D10").Value = 1123800
$ws.Range("E10").Value = 1233100
$ws.Range("F10").Value = 1122100
$ws.Range("G10").Value = 1351700
$ws.Range("H10").Value = 1277700
$ws.Range("I10").Value = 1311400
$ws.Range("J10").Value = 1280900

$ws.Range("D14").Value = 1100
$ws.Range("E14").Value = 19000
$ws.Range("F14").Value = -148400
$ws.Range("I14").Value = 3400

$ws.Range("D17").Value = 2109500
$ws.Range("E17").Value = 2043800
$ws.Range("F17").Value = 1918600
$ws.Range("G17").Value = 2118000
$ws.Range("H17").Value = 2057000
$ws.Range("I17").Value = 2299000
$ws.Range("J17").Value = 1940100

$ws.Range("D18").Value = 307600
$ws.Range("E18").Value = 455800
$ws.Range("F18").Value = 273800
$ws.Range("G18").Value = 412000
$ws.Range("H18").Value = 302800
$ws.Range("I18").Value = 405700
$ws.Range("J18").Value = 363300

$ws.Range("D20").Value = -117100
$ws.Range("E20").Value = -89200
$ws.Range("F20").Value = -112400
$ws.Range("G20").Value = -98900
$ws.Range("H20").Value = -77500
$ws.Range("I20").Value = -107800
$ws.Range("J20").Value = -95700

$ws.Range("D21").Value = 538500
$ws.Range("E21").Value = 357200
$ws.Range("F21").Value = 593500
$ws.Range("G21").Value = 349800
$ws.Range("H21").Value = 632600
$ws.Range("I21").Value = 323400
$ws.Range("J21").Value = 639500

$ws.Range("D23").Value = 190600
$ws.Range("E23").Value = 366600
$ws.Range("F23").Value = 161400
$ws.Range("G23").Value = 313100
$ws.Range("H23").Value = 225400
$ws.Range("I23").Value = 298000
$ws.Range("J23").Value = 267600

$ws.Range("D24").Value = 55700
$ws.Range("E24").Value = 108300
$ws.Range("F24").Value = 26900
$ws.Range("G24").Value = 42000
$ws.Range("I24").Value = 30300
$ws.Range("J24").Value = 26600

$ws.Range("D26").Value = 134900
$ws.Range("E26").Value = 258300
$ws.Range("F26").Value = 134500
$ws.Range("G26").Value = 271100
$ws.Range("H26").Value = 217100
$ws.Range("I26").Value = 267600
$ws.Range("J26").Value = 241000

$ws.Range("D27").Value = 23600
$ws.Range("E27").Value = 103800
$ws.Range("G27").Value = 150700
$ws.Range("H27").Value = 110600
$ws.Range("I27").Value = 156100
$ws.Range("J27").Value = 136300

$ws.Range("E29").Value = 8800
$ws.Range("F29").Value = 136800

$ws.Range("D32").Value = 117100
$ws.Range("E32").Value = 89200
$ws.Range("F32").Value = 112400
$ws.Range("G32").Value = 98900
$ws.Range("H32").Value = 77500
$ws.Range("I32").Value = 107800
$ws.Range("J32").Value = 95700

$ws.Range("D33").Value = 23600
$ws.Range("E33").Value = 112600
$ws.Range("F33").Value = 147000
$ws.Range("G33").Value = 150700
$ws.Range("H33").Value = 110600
$ws.Range("I33").Value = 156100
$ws.Range("J33").Value = 136300

$ws.Range("D35").Value = 23600
$ws.Range("E35").Value = 112600
$ws.Range("F35").Value = 147000
$ws.Range("G35").Value = 150700
$ws.Range("H35").Value = 110600
$ws.Range("I35").Value = 156100
$ws.Range("J35").Value = 136300

$ws.Range("D41").Value = 1073400
$ws.Range("E41").Value = 1482600
$ws.Range("F41").Value = 1286400
$ws.Range("G41").Value = 605200
$ws.Range("H41").Value = 718200
$ws.Range("I41").Value = 955800
$ws.Range("J41").Value = 836600

$ws.Range("D42").Value = 125600
$ws.Range("E42").Value = 217600
$ws.Range("F42").Value = 405600
$ws.Range("G42").Value = 57700
$ws.Range("J42").Value = 44600

$ws.Range("D43").Value = 992800
$ws.Range("E43").Value = 873800
$ws.Range("F43").Value = 414300
$ws.Range("G43").Value = 495800
$ws.Range("H43").Value = 563400
$ws.Range("I43").Value = 519200
$ws.Range("J43").Value = 628200

$ws.Range("D44").Value = 152000
$ws.Range("E44").Value = 116100
$ws.Range("F44").Value = 137200
$ws.Range("G44").Value = 120100
$ws.Range("H44").Value = 148200
$ws.Range("I44").Value = 98600
$ws.Range("J44").Value = 128200

$ws.Range("D45").Value = 582200
$ws.Range("E45").Value = 651000
$ws.Range("F45").Value = 1258600
$ws.Range("G45").Value = 1334400
$ws.Range("H45").Value = 1179600
$ws.Range("I45").Value = 991700
$ws.Range("J45").Value = 932400

$ws.Range("D46").Value = 2925900
$ws.Range("E46").Value = 3340900
$ws.Range("F46").Value = 3502100
$ws.Range("G46").Value = 2613300
$ws.Range("H46").Value = 2610100
$ws.Range("I46").Value = 2565500
$ws.Range("J46").Value = 2569800

$ws.Range("D47").Value = 448300
$ws.Range("E47").Value = 469400
$ws.Range("F47").Value = 417800
$ws.Range("G47").Value = 307000
$ws.Range("H47").Value = 281300
$ws.Range("I47").Value = 243200
$ws.Range("J47").Value = 250300

$ws.Range("D48").Value = 4291600
$ws.Range("E48").Value = 3802600
$ws.Range("F48").Value = 3193800
$ws.Range("G48").Value = 3037100
$ws.Range("H48").Value = 2931800
$ws.Range("I48").Value = 2757700
$ws.Range("J48").Value = 2632300

$ws.Range("D49").Value = 3563100
$ws.Range("E49").Value = 3532700
$ws.Range("F49").Value = 3775700
$ws.Range("G49").Value = 3831500
$ws.Range("H49").Value = 3614100
$ws.Range("I49").Value = 3657200
$ws.Range("J49").Value = 3667700

$ws.Range("D52").Value = 651100
$ws.Range("E52").Value = 641000
$ws.Range("F52").Value = 366100
$ws.Range("G52").Value = 349300
$ws.Range("H52").Value = 343700
$ws.Range("I52").Value = 299700
$ws.Range("J52").Value = 315400

$ws.Range("D54").Value = 11880000
$ws.Range("E54").Value = 11786600
$ws.Range("F54").Value = 11255500
$ws.Range("G54").Value = 10138200
$ws.Range("H54").Value = 9780900
$ws.Range("I54").Value = 9523400
$ws.Range("J54").Value = 9435500

$ws.Range("D57").Value = 238900
$ws.Range("E57").Value = 266000
$ws.Range("F57").Value = 245700
$ws.Range("G57").Value = 347900
$ws.Range("H57").Value = 326600
$ws.Range("I57").Value = 317700
$ws.Range("J57").Value = 371100

$ws.Range("D58").Value = 1214200
$ws.Range("E58").Value = 263400
$ws.Range("F58").Value = 67400
$ws.Range("G58").Value = 58200
$ws.Range("H58").Value = 296700
$ws.Range("I58").Value = 494100
$ws.Range("J58").Value = 1409600

$ws.Range("D59").Value = 1464300
$ws.Range("E59").Value = 1397600
$ws.Range("F59").Value = 1498500
$ws.Range("G59").Value = 1388000
$ws.Range("H59").Value = 1450700
$ws.Range("I59").Value = 1416600
$ws.Range("J59").Value = 1365400

$ws.Range("D60").Value = 2917400
$ws.Range("E60").Value = 1927000
$ws.Range("F60").Value = 1811600
$ws.Range("G60").Value = 1794200
$ws.Range("H60").Value = 2074000
$ws.Range("I60").Value = 2228400
$ws.Range("J60").Value = 3146000

$ws.Range("D61").Value = 5521000
$ws.Range("E61").Value = 6193400
$ws.Range("F61").Value = 5843100
$ws.Range("G61").Value = 5749200
$ws.Range("H61").Value = 5248900
$ws.Range("I61").Value = 4852300
$ws.Range("J61").Value = 3968700

$ws.Range("D62").Value = 900600
$ws.Range("E62").Value = 1053300
$ws.Range("F62").Value = 766500
$ws.Range("G62").Value = 712500
$ws.Range("H62").Value = 673000
$ws.Range("I62").Value = 743100
$ws.Range("J62").Value = 715800

$ws.Range("D66").Value = 9642100
$ws.Range("E66").Value = 9342000
$ws.Range("F66").Value = 8792300
$ws.Range("G66").Value = 8606200
$ws.Range("H66").Value = 8318300
$ws.Range("I66").Value = 8119100
$ws.Range("J66").Value = 8119300

$ws.Range("D72").Value = 667400
$ws.Range("E72").Value = 852000
$ws.Range("F72").Value = 881000
$ws.Range("G72").Value = -69200
$ws.Range("H72").Value = -140400
$ws.Range("I72").Value = -87500
$ws.Range("J72").Value = -168800

$ws.Range("D76").Value = 2237900
$ws.Range("E76").Value = 2444600
$ws.Range("F76").Value = 2463200
$ws.Range("G76").Value = 1532000
$ws.Range("H76").Value = 1462600
$ws.Range("I76").Value = 1404300
$ws.Range("J76").Value = 1316200

$ws.Range("D89").Value = 284200
$ws.Range("E89").Value = 864200
$ws.Range("F89").Value = 647000
$ws.Range("G89").Value = 829400
$ws.Range("H89").Value = 438500
$ws.Range("I89").Value = 950600
$ws.Range("J89").Value = 644100

$ws.Range("D94").Value = -254500
$ws.Range("E94").Value = -566500
$ws.Range("F94").Value = -660400
$ws.Range("G94").Value = -723800
$ws.Range("H94").Value = -557700
$ws.Range("I94").Value = -591000
$ws.Range("J94").Value = -707300

$ws.Range("D100").Value = -442300
$ws.Range("E100").Value = -93200
$ws.Range("F100").Value = 689100
$ws.Range("G100").Value = -214900
$ws.Range("H100").Value = -118600
$ws.Range("I100").Value = -241800
$ws.Range("J100").Value = -107600

$ws.Range("D101").Value = 3400

$ws.Range("D102").Value = -409200
$ws.Range("E102").Value = 196200
$ws.Range("F102").Value = 680400
$ws.Range("G102").Value = -112200
$ws.Range("H102").Value = -237600
$ws.Range("I102").Value = 119200
$ws.Range("J102").Value = -175300
